# fix: prevent hidden columns from being labeled upon detecting changes
#
# Rows that are the first row of a "Code" group but whose only detected
# "change" came from a column that should have been skipped now lose the
# erroneous "AENDERUNG" (change) label in column L, and - where that first
# row is itself such a false positive - the whole row's formatting is
# restored to the plain "group header" look (same look already used by the
# correctly-behaving groups in row 2 and row 9).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows where the group header row itself was wrongly flagged: restore the
# "no actual change" group-header formatting (copied from the known-good
# row 2 template) across the whole row, then deal with column L below.
$HEADER_ROWS = @(13, 17, 23, 27, 34, 40, 63, 67, 141)

foreach ($r in $HEADER_ROWS) {
    $ws.Range("A2:K2").Copy()
    $ws.Range("A" + $r + ":K" + $r).PasteSpecial(-4122)

    $ws.Range("M2:V2").Copy()
    $ws.Range("M" + $r + ":V" + $r).PasteSpecial(-4122)
}

# All rows whose column L carried a false-positive "AENDERUNG" label that
# must be removed (this set includes the header rows above as well as
# plain detail rows where only a hidden column looked changed).
$ALL_L_ROWS = @(13, 14, 15, 16, 17, 18, 19, 20, 21, 22, 23, 24, 25, 26, 27, 28, 29, 30, 31, 32, 33, 34, 35, 36, 38, 39, 40, 41, 42, 43, 63, 67, 135, 141, 142, 144, 147, 148, 149, 152, 153, 154, 155, 157, 158, 159, 160, 162, 164, 165, 166, 168, 169, 170, 171, 173, 174, 175, 177, 178, 179)

foreach ($r in $ALL_L_ROWS) {
    $ws.Range("L2").Copy()
    $ws.Range("L" + $r).PasteSpecial(-4122)
    $ws.Range("L" + $r).ClearContents()
}

$excel.CutCopyMode = 0
